$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.296.51"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.974.60"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.75"
$ws.Range("E5").Value = "  -4.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.46"
$ws.Range("E7").Value = "  -11.94%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -6.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.53"
$ws.Range("E10").Value = "  -6.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0832"
$ws.Range("E11").Value = "  +4.35%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.89"
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.854"
$ws.Range("E14").Value = "  -9.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.264.70"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("E16").Value = "  -8.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.39"
$ws.Range("E17").Value = "  -6.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.997.09"
$ws.Range("E18").Value = "  -2.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.200.24"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.22"
$ws.Range("E20").Value = "  -4.95%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0877"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.23"
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.91"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").Value = "  -5.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.75"
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.13"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.69"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.83"
$ws.Range("E33").Value = "  -7.62%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0667"
$ws.Range("E34").Value = "  +4.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.38"
$ws.Range("E35").Value = "  -7.18%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  -9.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("E40").Value = "  -7.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.21"
$ws.Range("E41").Value = "  -5.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0956"
$ws.Range("E42").Value = "  -8.19%  "
$ws.Range("E43").Value = "  -6.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0212"
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("E45").Value = "  -7.50%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.98"
$ws.Range("E46").Value = "  -11.26%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.10"
$ws.Range("E47").Value = "  -6.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.357.41"
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.37"
$ws.Range("E49").Value = "  -8.41%  "
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.56"
$ws.Range("E51").Value = "  -9.32%  "
